$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add payment amounts for the existing rows 27-30 (rightmost "Оплата" column)
$ws.Range("E27").Value() = 800
$ws.Range("E28").Value() = 800
$ws.Range("E29").Value() = 800
$ws.Range("E30").Value() = 800

# Append a new lesson row (#34 - "Insert sort") with its config-file note
$ws.Range("A31").Value() = 34
$ws.Range("B31").Value() = 44572
$ws.Range("C31").Value() = "Insert sort"
$ws.Range("D31").Value() = "D:\Teaching\12.Algorithms\Tasks\Tasks"

# Match the saved selection/cursor position from the source workbook
$ws.Range("E31").Select() | Out-Null
